# Update the "想去人数" (want-to-go count) values in column F for a set of
# events on sheets "展览" and "全部类型". These two sheets both list the same
# events, so the same updates are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 519
    7  = 92
    10 = 6534
    11 = 222
    12 = 355
    13 = 2714
    14 = 166
    15 = 283
    17 = 513
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates4 = @{
    5  = 519
    9  = 92
    13 = 6534
    15 = 222
    16 = 355
    17 = 2714
    18 = 166
    19 = 283
    21 = 513
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
